# Supervision section reformatted (data/supervision_pos_es.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new row at position 5 ---
# In the previous layout, the 2nd supervision entry (Francisco Javier Flores)
# had no dedicated "note" row, so the 3rd/4th entries started one row earlier.
# The new layout gives every entry its own note row, so we insert one row to
# push everything from the old row 5 onward down by one (18 -> 20 rows total).
$ws.Rows("5:5").Insert()

# --- Row 1: headers (unchanged) ---
$ws.Range('A1').Value = 'what'
$ws.Range('B1').Value = 'when'
$ws.Range('C1').Value = 'with'
$ws.Range('D1').Value = 'where'
$ws.Range('E1').Value = 'why'

# --- Entry 1: PhD in Neuroscience (Milena Vasquez-Amezquita) ---
$ws.Range('A2').Value = 'PhD in Neuroscience'
$ws.Range('B2').Value = '2015 - 2018'
$ws.Range('C2').Value = '\href{https://www.researchgate.net/profile/Milena-Vasquez-Amezquita}{Milena Vásquez-Amézquita}'
$ws.Range('D2').Value = '\href{https://www.uv.es/}{Universitat de València}, España'
$ws.Range('E2').Value = 'Tésis \textbf{\textit{(Summa Cum Laude)}}: \textit{\href{http://hdl.handle.net/10550/67639}{Preferencias sexuales típicas y atípicas según sexo y edad de los estímulos: Utilidad de la técnica de rastreo ocular} [Typical and atypical sexual preferences according to sex and age of the stimuli: Usefulness of the eye tracking technique]}'
$ws.Range('A3').Value = $null
$ws.Range('B3').Value = $null
$ws.Range('C3').Value = $null
$ws.Range('D3').Value = $null
$ws.Range('E3').Value = 'Supervisión conjunta con  Alicia Salvador'

# --- Entry 2: Professional Doctorate in Counselling Psychology (Francisco Javier Flores) ---
$ws.Range('A4').Value = 'Professional Doctorate in Counselling Psychology'
$ws.Range('B4').Value = '2015 - 2018'
$ws.Range('C4').Value = '\href{https://www.researchgate.net/profile/Francisco-Flores-14}{Francisco Javier Flores}'
$ws.Range('D4').Value = '\href{https://www.uel.ac.uk/}{U. of East London}, Reino Unido'
$ws.Range('E4').Value = 'Tésis: \textit{ What sense do people make of the functions of their ’behaviours that may be causing problems in their everyday life’? A hybrid deductive/inductive template analysis}'
$ws.Range('A5').Value = $null
$ws.Range('B5').Value = $null
$ws.Range('C5').Value = $null
$ws.Range('D5').Value = $null
$ws.Range('E5').Value = 'Supervisión conjunta con Lisa Chiara Fellin'

# --- Entry 3: Psychological Research Methods (Evolutionary Psychology) MSc (Julia Sanz-Vidania) ---
$ws.Range('A6').Value = 'Psychological Research Methods (Evolutionary Psychology) MSc'
$ws.Range('B6').Value = '2013 - 2014'
$ws.Range('C6').Value = 'Julia Sanz-Vidania'
$ws.Range('D6').Value = '\href{https://www.stir.ac.uk/}{University of Stirling}, Reino Unido'
$ws.Range('E6').Value = 'Trabajo de grado \textbf{\textit{(Meritorio)}}: \textit{Sexy Chat: Mate-Choice Preferences for Speech Content in the Absence of Auditory Cues}'
$ws.Range('A7').Value = $null
$ws.Range('B7').Value = $null
$ws.Range('C7').Value = $null
$ws.Range('D7').Value = $null
$ws.Range('E7').Value = 'Supervisión conjunta con \href{https://www.scraigroberts.com/}{S Craig Roberts}'

# --- Entry 4: Maestría en Psicología (Adrián Acosta Guerrero) ---
$ws.Range('A8').Value = 'Maestría en Psicología'
$ws.Range('B8').Value = '2019 - 2020'
$ws.Range('C8').Value = 'Adrián Acosta Guerrero'
$ws.Range('D8').Value = '\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia'
$ws.Range('E8').Value = 'Trabajo de grado \textbf{\textit{(Meritorio)}}: \textit{\href{http://hdl.handle.net/20.500.12495/4416}{La voz como predictor de sintomatología asociada a depresión y ansiedad} [Voice as a predictor of symptomatology associated with depression and anxiety]}'
$ws.Range('A9').Value = $null
$ws.Range('B9').Value = $null
$ws.Range('C9').Value = $null
$ws.Range('D9').Value = $null
$ws.Range('E9').Value = 'Supervisión conjunta con \href{https://www.researchgate.net/profile/Milena-Vasquez-Amezquita}{Milena Vásquez-Amézquita}'

# --- Row heights: entry rows are taller to fit wrapped text ---
$ws.Rows("2:2").RowHeight = 45
$ws.Rows("4:4").RowHeight = 30
$ws.Rows("6:6").RowHeight = 30
$ws.Range("A8:E8").RowHeight = 31.5

# --- Column E width adjustment ---
$ws.Columns("E:E").ColumnWidth = 118.42578125

# --- Sheet view: scrolled right, selection on E8 ---
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("E8").Select()
